$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "I AM One Who A B" sample grid and lay out the new
# 4-column team roster header row.
$ws.Cells.Clear()

$ws.Range("A1").Value = "Team Name"
$ws.Range("B1").Value = "Team Email"
$ws.Range("C1").Value = "Team Phone"
$ws.Range("D1").Value = "Team Budget"

# Style the header row: big white text on a near-black fill (order
# matters - fill first, then font, keeps the style table minimal).
$headerRange = $ws.Range("A1:D1")
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 1
$headerRange.Font.Size = 22
$headerRange.Font.ThemeColor = 2

# Taller header row to match the larger font.
$ws.Rows.Item(1).RowHeight = 28.5

# Wider columns for the new headers.
$ws.Columns.Item("A:B").ColumnWidth = 26.5703125
$ws.Columns.Item("C:D").ColumnWidth = 25.42578125

# Leave the selection where the user ended up after entering data.
$ws.Range("B19").Select()

# Touch the page setup (portrait) so it is persisted with the sheet.
$ws.PageSetup.Orientation = 1
